$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 152, pushing existing rows 152..214 down to 153..215.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row 152 with its data (same "template" columns as
# the surrounding Oregano / Lo Valledor rows, new Fecha/Volumen/Precios values).
$ws.Cells.Item(152, 1).Value = 6
$ws.Cells.Item(152, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(152, 3).Value = "Metropolitana"
$ws.Cells.Item(152, 4).Value = 44726
$ws.Cells.Item(152, 5).Value = 13
$ws.Cells.Item(152, 6).Value = 100112029
$ws.Cells.Item(152, 7).Value = "Orégano"
$ws.Cells.Item(152, 8).Value = "Sin especificar"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 46
$ws.Cells.Item(152, 11).Value = 12000
$ws.Cells.Item(152, 12).Value = 13000
$ws.Cells.Item(152, 13).Value = 12478
$ws.Cells.Item(152, 14).Value = "`$/docena de atados"
$ws.Cells.Item(152, 15).Value = "Región Metropolitana"
$ws.Cells.Item(152, 16).Value = 4159
$ws.Cells.Item(152, 17).Value = 3
$ws.Cells.Item(152, 18).Value = "Hortaliza"
